$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3

$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 0

$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 1

$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 2

$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 1

$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1
